$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 688.4286
$ws.Range("I6").Value = 60.3
$ws.Range("J6").Value = 2258.75
$ws.Range("K6").Value = 180.9
$ws.Range("L6").Value = 6776.25
$ws.Range("M6").Value = -68.89999999999998
$ws.Range("N6").Value = -7000.25

# ALC row 51
$ws.Range("H51").Value = 2979.6667
$ws.Range("I51").Value = 2969
$ws.Range("K51").Value = 2969
$ws.Range("M51").Value = -2485

# ALC row 70
$ws.Range("H70").Value = 5428.091
$ws.Range("I70").Value = 1253
$ws.Range("J70").Value = 8318.538
$ws.Range("K70").Value = 3759
$ws.Range("L70").Value = 24955.614
$ws.Range("M70").Value = -3489
$ws.Range("N70").Value = -25495.614

# ALC row 73
$ws.Range("H73").Value = 5428.091
$ws.Range("I73").Value = 1253
$ws.Range("J73").Value = 8318.538
$ws.Range("K73").Value = 3759
$ws.Range("L73").Value = 24955.614
$ws.Range("M73").Value = -2823
$ws.Range("N73").Value = -26827.614

# ALC row 87
$ws.Range("H87").Value = 75838
$ws.Range("J87").Value = 75838
$ws.Range("L87").Value = 75838
$ws.Range("N87").Value = -78334

# ALC row 90
$ws.Range("H90").Value = 75838
$ws.Range("J90").Value = 75838
$ws.Range("L90").Value = 227514
$ws.Range("N90").Value = -239994

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# ARM row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# ARM row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# ARM row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# ARM row 122
$ws.Range("H122").Value = 1125.5
$ws.Range("I122").Value = 995.6667
$ws.Range("K122").Value = 2987.0001
$ws.Range("M122").Value = -537.0001000000002

# ARM row 132
$ws.Range("H132").Value = 1233.3334
$ws.Range("I132").Value = 1233.3334
$ws.Range("K132").Value = 3700.0002
$ws.Range("M132").Value = -1170.0002

# ARM row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 2668.4
$ws.Range("I86").Value = 1568.5416
$ws.Range("J86").Value = 7067.8335
$ws.Range("K86").Value = 1568.5416
$ws.Range("L86").Value = 7067.8335
$ws.Range("M86").Value = -445.5416
$ws.Range("N86").Value = -9313.833500000001

# BSM row 89
$ws.Range("H89").Value = 2668.4
$ws.Range("I89").Value = 1568.5416
$ws.Range("J89").Value = 7067.8335
$ws.Range("K89").Value = 7842.708000000001
$ws.Range("L89").Value = 35339.1675
$ws.Range("M89").Value = -2226.708000000001
$ws.Range("N89").Value = -46571.1675

# BSM row 134
$ws.Range("H134").Value = 2940
$ws.Range("I134").Value = 880
$ws.Range("K134").Value = 2640
$ws.Range("M134").Value = -105

$ws = $wb.Worksheets.Item("CRP")
# CRP row 25
$ws.Range("H25").Value = 4450
$ws.Range("I25").Value = 4450
$ws.Range("K25").Value = 4450
$ws.Range("M25").Value = -4276

# CRP row 141
$ws.Range("H141").Value = 82805.89999999999
$ws.Range("J141").Value = 82805.89999999999
$ws.Range("L141").Value = 82805.89999999999
$ws.Range("N141").Value = -93165.89999999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 1599.2222
$ws.Range("I5").Value = 984.7143
$ws.Range("J5").Value = 3750
$ws.Range("K5").Value = 2954.1429
$ws.Range("L5").Value = 11250
$ws.Range("M5").Value = -2842.1429
$ws.Range("N5").Value = -11474

# CUL row 132
$ws.Range("H132").Value = 1466.6666
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 21600
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -26660

# CUL row 134
$ws.Range("H134").Value = 4447.25
$ws.Range("I134").Value = 4447.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13341.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8271.75
$ws.Range("N134").ClearContents()

# CUL row 135
$ws.Range("H135").Value = 1599.2222
$ws.Range("I135").Value = 984.7143
$ws.Range("J135").Value = 3750
$ws.Range("K135").Value = 8862.4287
$ws.Range("L135").Value = 33750
$ws.Range("M135").Value = -6327.4287
$ws.Range("N135").Value = -38820

# CUL row 136
$ws.Range("H136").Value = 8515
$ws.Range("I136").Value = 6030
$ws.Range("K136").Value = 18090
$ws.Range("M136").Value = -12990

# CUL row 137
$ws.Range("H137").Value = 5150
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 5300
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 15900
$ws.Range("M137").Value = -9900
$ws.Range("N137").Value = -26100

$ws = $wb.Worksheets.Item("GSM")
# GSM row 55
$ws.Range("H55").Value = 7257.25
$ws.Range("I55").Value = 6676.6665
$ws.Range("J55").Value = 8999
$ws.Range("K55").Value = 6676.6665
$ws.Range("L55").Value = 8999
$ws.Range("M55").Value = -6349.6665
$ws.Range("N55").Value = -9653

# GSM row 107
$ws.Range("H107").Value = 214.44444
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 306
$ws.Range("K107").Value = 100
$ws.Range("L107").Value = 306
$ws.Range("M107").Value = 1820
$ws.Range("N107").Value = -4146

# GSM row 132
$ws.Range("H132").Value = 4415
$ws.Range("I132").Value = 4693.222
$ws.Range("K132").Value = 14079.666
$ws.Range("M132").Value = -11549.666

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 4131.4614
$ws.Range("I40").Value = 4131.4614
$ws.Range("K40").Value = 4131.4614
$ws.Range("M40").Value = -3995.4614

# LTW row 61
$ws.Range("H61").Value = 4385.857
$ws.Range("I61").Value = 2425.5
$ws.Range("K61").Value = 2425.5
$ws.Range("M61").Value = -2223.5

# LTW row 68
$ws.Range("H68").Value = 8731.25
$ws.Range("J68").Value = 9666.666999999999
$ws.Range("L68").Value = 9666.666999999999
$ws.Range("N68").Value = -11164.667

# LTW row 71
$ws.Range("H71").Value = 8731.25
$ws.Range("J71").Value = 9666.666999999999
$ws.Range("L71").Value = 48333.335
$ws.Range("N71").Value = -55821.335

# LTW row 113
$ws.Range("H113").Value = 4385.857
$ws.Range("I113").Value = 2425.5
$ws.Range("K113").Value = 2425.5
$ws.Range("M113").Value = -255.5

# LTW row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2
$ws.Range("H2").Value = 132.14285
$ws.Range("I2").Value = 132.14285
$ws.Range("K2").Value = 132.14285
$ws.Range("M2").Value = -20.14285000000001

# WVR row 4
$ws.Range("H4").Value = 1408.3334
$ws.Range("I4").Value = 1408.3334
$ws.Range("K4").Value = 1408.3334
$ws.Range("M4").Value = -1295.3334

# WVR row 5
$ws.Range("H5").Value = 13714286
$ws.Range("I5").Value = 15000000
$ws.Range("J5").Value = 10500000
$ws.Range("K5").Value = 15000000
$ws.Range("L5").Value = 10500000
$ws.Range("M5").Value = -14999888
$ws.Range("N5").Value = -10500224

# WVR row 95
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492

# WVR row 96
$ws.Range("H96").Value = 1298.5
$ws.Range("I96").Value = 1281.75
$ws.Range("J96").Value = 1320.8334
$ws.Range("K96").Value = 1281.75
$ws.Range("L96").Value = 1320.8334
$ws.Range("M96").Value = 91.25
$ws.Range("N96").Value = -4066.8334

# WVR row 132
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
